$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45982
$ws.Range("B2").Value = 73.18000000000001
$ws.Range("C2").Value = 64.36
$ws.Range("D2").Value = 62.35
$ws.Range("E2").Value = 58.59
$ws.Range("F2").Value = 58.28
$ws.Range("G2").Value = 52.98
$ws.Range("H2").Value = 59.5
$ws.Range("I2").Value = 71.94
$ws.Range("J2").Value = 73.86
$ws.Range("K2").Value = 27.83
$ws.Range("L2").Value = 0.31
$ws.Range("M2").Value = 0.8
$ws.Range("N2").Value = 6.44
$ws.Range("O2").Value = 7.08
$ws.Range("P2").Value = 21.56
$ws.Range("Q2").Value = 43.47
$ws.Range("R2").Value = 60.43
$ws.Range("S2").Value = 90.11
$ws.Range("T2").Value = 113.38
$ws.Range("U2").Value = 147.43
$ws.Range("V2").Value = 135.34
$ws.Range("W2").Value = 104.31
$ws.Range("X2").Value = 99.19
$ws.Range("Y2").Value = 92.27
$ws.Range("Z2").Value = 63.54
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 107.78
$ws.Range("AD2").Value = 130.4
$ws.Range("AF2").Value = 119.82
$ws.Range("AG2").Value = "2h-16h"
